# Updates Sheets/Odin_Profits.xlsx data per scheduled runner refresh
$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(70, 8).Value = 10100.5  # H70
$ws.Cells.Item(70, 9).Value = 1250  # I70
$ws.Cells.Item(70, 10).Value = 14525.75  # J70
$ws.Cells.Item(70, 11).Value = 3750  # K70
$ws.Cells.Item(70, 12).Value = 43577.25  # L70
$ws.Cells.Item(70, 13).Value = -3480  # M70
$ws.Cells.Item(70, 14).Value = -44117.25  # N70
$ws.Cells.Item(73, 8).Value = 10100.5  # H73
$ws.Cells.Item(73, 9).Value = 1250  # I73
$ws.Cells.Item(73, 10).Value = 14525.75  # J73
$ws.Cells.Item(73, 11).Value = 3750  # K73
$ws.Cells.Item(73, 12).Value = 43577.25  # L73
$ws.Cells.Item(73, 13).Value = -2814  # M73
$ws.Cells.Item(73, 14).Value = -45449.25  # N73
$ws.Cells.Item(98, 8).Value = 1186.8572  # H98
$ws.Cells.Item(98, 9).Value = 717.5  # I98
$ws.Cells.Item(98, 11).Value = 717.5  # K98
$ws.Cells.Item(98, 13).Value = 780.5  # M98
$ws.Cells.Item(103, 8).Value = 933  # H103
$ws.Cells.Item(103, 10).Value = 900  # J103
$ws.Cells.Item(103, 12).Value = 2700  # L103
$ws.Cells.Item(103, 14).Value = -3872  # N103
$ws.Cells.Item(112, 8).Value = 2418.2  # H112
$ws.Cells.Item(112, 10).Value = 2208.2273  # J112
$ws.Cells.Item(112, 12).Value = 6624.6819  # L112
$ws.Cells.Item(112, 14).Value = -8840.6819  # N112
$ws.Cells.Item(116, 8).Value = 6948538  # H116
$ws.Cells.Item(116, 9).Value = 10104610  # I116
$ws.Cells.Item(116, 10).Value = 5180.8  # J116
$ws.Cells.Item(116, 11).Value = 10104610  # K116
$ws.Cells.Item(116, 12).Value = 5180.8  # L116
$ws.Cells.Item(116, 13).Value = -10101168  # M116
$ws.Cells.Item(116, 14).Value = -12064.8  # N116
$ws.Cells.Item(122, 8).Value = 1186.8572  # H122
$ws.Cells.Item(122, 9).Value = 717.5  # I122
$ws.Cells.Item(122, 11).Value = 2152.5  # K122
$ws.Cells.Item(122, 13).Value = 297.5  # M122
$ws.Cells.Item(138, 8).Value = 4232.054  # H138
$ws.Cells.Item(138, 9).Value = 2120.8823  # I138
$ws.Cells.Item(138, 11).Value = 6362.646900000001  # K138
$ws.Cells.Item(138, 13).Value = -1222.646900000001  # M138

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 3763.2812  # H2
$ws.Cells.Item(2, 9).Value = 1286.4166  # I2
$ws.Cells.Item(2, 10).Value = 11193.875  # J2
$ws.Cells.Item(2, 11).Value = 1286.4166  # K2
$ws.Cells.Item(2, 12).Value = 11193.875  # L2
$ws.Cells.Item(2, 13).Value = -1173.4166  # M2
$ws.Cells.Item(2, 14).Value = -11419.875  # N2
$ws.Cells.Item(32, 8).Value = 4835350  # H32
$ws.Cells.Item(32, 9).Value = 2913.2222  # I32
$ws.Cells.Item(32, 10).Value = 22232122  # J32
$ws.Cells.Item(32, 11).Value = 2913.2222  # K32
$ws.Cells.Item(32, 12).Value = 22232122  # L32
$ws.Cells.Item(32, 13).Value = -2626.2222  # M32
$ws.Cells.Item(32, 14).Value = -22232696  # N32
$ws.Cells.Item(61, 8).Value = 4499.213  # H61
$ws.Cells.Item(61, 9).Value = 4049.55  # I61
$ws.Cells.Item(61, 11).Value = 4049.55  # K61
$ws.Cells.Item(61, 13).Value = -3837.55  # M61
$ws.Cells.Item(116, 8).Value = 3763.2812  # H116
$ws.Cells.Item(116, 9).Value = 1286.4166  # I116
$ws.Cells.Item(116, 10).Value = 11193.875  # J116
$ws.Cells.Item(116, 11).Value = 1286.4166  # K116
$ws.Cells.Item(116, 12).Value = 11193.875  # L116
$ws.Cells.Item(116, 13).Value = 1007.5834  # M116
$ws.Cells.Item(116, 14).Value = -15781.875  # N116
$ws.Cells.Item(122, 8).Value = 0  # H122
$ws.Cells.Item(122, 9).Value = 0  # I122
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 11).Value = 0  # K122
$ws.Cells.Item(122, 12).ClearContents()  # L122
$ws.Cells.Item(122, 13).ClearContents()  # M122
$ws.Cells.Item(122, 14).Value = 0  # N122
$ws.Cells.Item(136, 8).Value = 4499.213  # H136
$ws.Cells.Item(136, 9).Value = 4049.55  # I136
$ws.Cells.Item(136, 11).Value = 12148.65  # K136
$ws.Cells.Item(136, 13).Value = -9598.650000000001  # M136

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 3763.2812  # H3
$ws.Cells.Item(3, 9).Value = 1286.4166  # I3
$ws.Cells.Item(3, 10).Value = 11193.875  # J3
$ws.Cells.Item(3, 11).Value = 1286.4166  # K3
$ws.Cells.Item(3, 12).Value = 11193.875  # L3
$ws.Cells.Item(3, 13).Value = -1172.4166  # M3
$ws.Cells.Item(3, 14).Value = -11421.875  # N3
$ws.Cells.Item(60, 8).Value = 90000  # H60
$ws.Cells.Item(60, 10).Value = 90000  # J60
$ws.Cells.Item(60, 12).Value = 90000  # L60
$ws.Cells.Item(60, 14).Value = -91198  # N60

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 23815234  # H16
$ws.Cells.Item(16, 9).Value = 41671252  # I16
$ws.Cells.Item(16, 10).Value = 7212.4443  # J16
$ws.Cells.Item(16, 11).Value = 41671252  # K16
$ws.Cells.Item(16, 12).Value = 7212.4443  # L16
$ws.Cells.Item(16, 13).Value = -41670965  # M16
$ws.Cells.Item(16, 14).Value = -7786.4443  # N16
$ws.Cells.Item(31, 8).Value = 5028.45  # H31
$ws.Cells.Item(31, 9).Value = 1316.875  # I31
$ws.Cells.Item(31, 10).Value = 7502.8335  # J31
$ws.Cells.Item(31, 11).Value = 1316.875  # K31
$ws.Cells.Item(31, 12).Value = 7502.8335  # L31
$ws.Cells.Item(31, 13).Value = -1021.875  # M31
$ws.Cells.Item(31, 14).Value = -8092.8335  # N31
$ws.Cells.Item(34, 8).Value = 5028.45  # H34
$ws.Cells.Item(34, 9).Value = 1316.875  # I34
$ws.Cells.Item(34, 10).Value = 7502.8335  # J34
$ws.Cells.Item(34, 11).Value = 1316.875  # K34
$ws.Cells.Item(34, 12).Value = 7502.8335  # L34
$ws.Cells.Item(34, 13).Value = -1114.875  # M34
$ws.Cells.Item(34, 14).Value = -7906.8335  # N34
$ws.Cells.Item(58, 8).Value = 28583286  # H58
$ws.Cells.Item(58, 9).Value = 41674010  # I58
$ws.Cells.Item(58, 10).Value = 21715.727  # J58
$ws.Cells.Item(58, 11).Value = 41674010  # K58
$ws.Cells.Item(58, 12).Value = 21715.727  # L58
$ws.Cells.Item(58, 13).Value = -41673807  # M58
$ws.Cells.Item(58, 14).Value = -22121.727  # N58
$ws.Cells.Item(94, 8).Value = 71432510  # H94
$ws.Cells.Item(94, 9).Value = 142858850  # I94
$ws.Cells.Item(94, 11).Value = 142858850  # K94
$ws.Cells.Item(94, 13).Value = -142858399  # M94
$ws.Cells.Item(113, 8).Value = 23815234  # H113
$ws.Cells.Item(113, 9).Value = 41671252  # I113
$ws.Cells.Item(113, 10).Value = 7212.4443  # J113
$ws.Cells.Item(113, 11).Value = 41671252  # K113
$ws.Cells.Item(113, 12).Value = 7212.4443  # L113
$ws.Cells.Item(113, 13).Value = -41669082  # M113
$ws.Cells.Item(113, 14).Value = -11552.4443  # N113
$ws.Cells.Item(134, 8).Value = 83344740  # H134
$ws.Cells.Item(134, 9).Value = 100005710  # I134
$ws.Cells.Item(134, 11).Value = 300017130  # K134
$ws.Cells.Item(134, 13).Value = -300014595  # M134
$ws.Cells.Item(136, 8).Value = 28583286  # H136
$ws.Cells.Item(136, 9).Value = 41674010  # I136
$ws.Cells.Item(136, 10).Value = 21715.727  # J136
$ws.Cells.Item(136, 11).Value = 125022030  # K136
$ws.Cells.Item(136, 12).Value = 65147.181  # L136
$ws.Cells.Item(136, 13).Value = -125019480  # M136
$ws.Cells.Item(136, 14).Value = -70247.181  # N136

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(37, 8).Value = 102555.5  # H37
$ws.Cells.Item(37, 10).Value = 102555.5  # J37
$ws.Cells.Item(37, 12).Value = 307666.5  # L37
$ws.Cells.Item(37, 14).Value = -307890.5  # N37
$ws.Cells.Item(106, 8).Value = 28300  # H106
$ws.Cells.Item(106, 10).Value = 28300  # J106
$ws.Cells.Item(106, 12).Value = 84900  # L106
$ws.Cells.Item(106, 14).Value = -86792  # N106

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value = 11854.4  # H113
$ws.Cells.Item(113, 9).Value = 6555  # I113
$ws.Cells.Item(113, 10).Value = 12443.223  # J113
$ws.Cells.Item(113, 11).Value = 6555  # K113
$ws.Cells.Item(113, 12).Value = 12443.223  # L113
$ws.Cells.Item(113, 13).Value = -4385  # M113
$ws.Cells.Item(113, 14).Value = -16783.223  # N113
$ws.Cells.Item(122, 8).Value = 5350.476  # H122
$ws.Cells.Item(122, 9).Value = 4316.273  # I122
$ws.Cells.Item(122, 10).Value = 6488.1  # J122
$ws.Cells.Item(122, 11).Value = 12948.819  # K122
$ws.Cells.Item(122, 12).Value = 19464.3  # L122
$ws.Cells.Item(122, 13).Value = -10498.819  # M122
$ws.Cells.Item(122, 14).Value = -24364.3  # N122
$ws.Cells.Item(126, 8).Value = 38474620  # H126
$ws.Cells.Item(126, 9).Value = 71433300  # I126
$ws.Cells.Item(126, 11).Value = 214299900  # K126
$ws.Cells.Item(126, 13).Value = -214297430  # M126
$ws.Cells.Item(132, 8).Value = 24393462  # H132
$ws.Cells.Item(132, 9).Value = 29414500  # I132
$ws.Cells.Item(132, 10).Value = 5571.2856  # J132
$ws.Cells.Item(132, 11).Value = 88243500  # K132
$ws.Cells.Item(132, 12).Value = 16713.8568  # L132
$ws.Cells.Item(132, 13).Value = -88240970  # M132
$ws.Cells.Item(132, 14).Value = -21773.8568  # N132
$ws.Cells.Item(137, 8).Value = 46990  # H137
$ws.Cells.Item(137, 10).Value = 46990  # J137
$ws.Cells.Item(137, 12).Value = 46990  # L137
$ws.Cells.Item(137, 14).Value = -57190  # N137

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 4098.8887  # H40
$ws.Cells.Item(40, 9).Value = 3480  # I40
$ws.Cells.Item(40, 10).Value = 4872.5  # J40
$ws.Cells.Item(40, 11).Value = 3480  # K40
$ws.Cells.Item(40, 12).Value = 4872.5  # L40
$ws.Cells.Item(40, 13).Value = -3344  # M40
$ws.Cells.Item(40, 14).Value = -5144.5  # N40
$ws.Cells.Item(122, 8).Value = 5699.3335  # H122
$ws.Cells.Item(122, 9).Value = 4799  # I122
$ws.Cells.Item(122, 11).Value = 14397  # K122
$ws.Cells.Item(122, 13).Value = -11947  # M122
$ws.Cells.Item(136, 8).Value = 38469730  # H136
$ws.Cells.Item(136, 9).Value = 83341750  # I136
$ws.Cells.Item(136, 10).Value = 7995.5713  # J136
$ws.Cells.Item(136, 11).Value = 250025250  # K136
$ws.Cells.Item(136, 12).Value = 23986.7139  # L136
$ws.Cells.Item(136, 13).Value = -250022700  # M136
$ws.Cells.Item(136, 14).Value = -29086.7139  # N136

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 20099.6  # H62
$ws.Cells.Item(62, 10).Value = 16937  # J62
$ws.Cells.Item(62, 12).Value = 16937  # L62
$ws.Cells.Item(62, 14).Value = -18185  # N62
$ws.Cells.Item(65, 8).Value = 20099.6  # H65
$ws.Cells.Item(65, 10).Value = 16937  # J65
$ws.Cells.Item(65, 12).Value = 84685  # L65
$ws.Cells.Item(65, 14).Value = -90925  # N65
$ws.Cells.Item(122, 8).Value = 15465.182  # H122
$ws.Cells.Item(122, 9).Value = 9188.786  # I122
$ws.Cells.Item(122, 10).Value = 26448.875  # J122
$ws.Cells.Item(122, 11).Value = 27566.358  # K122
$ws.Cells.Item(122, 12).Value = 79346.625  # L122
$ws.Cells.Item(122, 13).Value = -25116.358  # M122
$ws.Cells.Item(122, 14).Value = -84246.625  # N122
$ws.Cells.Item(136, 8).Value = 17870064  # H136
$ws.Cells.Item(136, 9).Value = 22738076  # I136
$ws.Cells.Item(136, 11).Value = 68214228  # K136
$ws.Cells.Item(136, 13).Value = -68211678  # M136
